$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.078459
$ws.Range("H2").Value = 0.235377
$ws.Range("I2").Value = 0.01931367871928443
$ws.Range("J2").Value = 0.01931367871928443
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.078459
$ws.Range("N2").Value = 0.235377
$ws.Range("O2").Value = 0.01931367871928443
$ws.Range("P2").Value = 0.01931367871928443
$ws.Range("Q2").Value = 0.006155814681
$ws.Range("R2").Value = 0.055402332129
$ws.Range("S2").Value = 0.0003730181856717404
$ws.Range("T2").Value = 0.0003730181856717404
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.078459
$ws.Range("H3").Value = 0.235377
$ws.Range("I3").Value = 0.01931367871928443
$ws.Range("J3").Value = 0.01931367871928443
$ws.Range("O3").Value = 0.5390230229402296
$ws.Range("P3").Value = 0.5390230229402296
$ws.Range("Q3").Value = 0.171801855371
$ws.Range("R3").Value = 1.546216698339
$ws.Range("S3").Value = 0.01041051748736508
$ws.Range("T3").Value = 0.01041051748736508
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.078459
$ws.Range("H4").Value = 0.235377
$ws.Range("I4").Value = 0.01931367871928443
$ws.Range("J4").Value = 0.01931367871928443
$ws.Range("M4").Value = 1.794192666666667
$ws.Range("N4").Value = 5.382578
$ws.Range("O4").Value = 0.441663298340486
$ws.Range("P4").Value = 0.441663298340486
$ws.Range("Q4").Value = 0.140770562434
$ws.Range("R4").Value = 1.266935061906
$ws.Range("S4").Value = 0.008530143046247616
$ws.Range("T4").Value = 0.008530143046247617
$ws.Range("I5").Value = 0.5390230229402296
$ws.Range("J5").Value = 0.5390230229402296
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.078459
$ws.Range("N5").Value = 0.235377
$ws.Range("O5").Value = 0.01931367871928443
$ws.Range("P5").Value = 0.01931367871928443
$ws.Range("Q5").Value = 0.171801855371
$ws.Range("R5").Value = 1.546216698339
$ws.Range("S5").Value = 0.01041051748736508
$ws.Range("T5").Value = 0.01041051748736508
$ws.Range("I6").Value = 0.5390230229402296
$ws.Range("J6").Value = 0.5390230229402296
$ws.Range("O6").Value = 0.5390230229402296
$ws.Range("P6").Value = 0.5390230229402296
$ws.Range("S6").Value = 0.2905458192596233
$ws.Range("T6").Value = 0.2905458192596233
$ws.Range("I7").Value = 0.5390230229402296
$ws.Range("J7").Value = 0.5390230229402296
$ws.Range("M7").Value = 1.794192666666667
$ws.Range("N7").Value = 5.382578
$ws.Range("O7").Value = 0.441663298340486
$ws.Range("P7").Value = 0.441663298340486
$ws.Range("Q7").Value = 3.928747868649555
$ws.Range("R7").Value = 35.35873081784599
$ws.Range("S7").Value = 0.2380666861932413
$ws.Range("T7").Value = 0.2380666861932413
$ws.Range("G8").Value = 1.794192666666667
$ws.Range("H8").Value = 5.382578
$ws.Range("I8").Value = 0.441663298340486
$ws.Range("J8").Value = 0.441663298340486
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.078459
$ws.Range("N8").Value = 0.235377
$ws.Range("O8").Value = 0.01931367871928443
$ws.Range("P8").Value = 0.01931367871928443
$ws.Range("Q8").Value = 0.140770562434
$ws.Range("R8").Value = 1.266935061906
$ws.Range("S8").Value = 0.008530143046247616
$ws.Range("T8").Value = 0.008530143046247617
$ws.Range("G9").Value = 1.794192666666667
$ws.Range("H9").Value = 5.382578
$ws.Range("I9").Value = 0.441663298340486
$ws.Range("J9").Value = 0.441663298340486
$ws.Range("O9").Value = 0.5390230229402296
$ws.Range("P9").Value = 0.5390230229402296
$ws.Range("Q9").Value = 3.928747868649555
$ws.Range("R9").Value = 35.35873081784599
$ws.Range("S9").Value = 0.2380666861932413
$ws.Range("T9").Value = 0.2380666861932413
$ws.Range("G10").Value = 1.794192666666667
$ws.Range("H10").Value = 5.382578
$ws.Range("I10").Value = 0.441663298340486
$ws.Range("J10").Value = 0.441663298340486
$ws.Range("M10").Value = 1.794192666666667
$ws.Range("N10").Value = 5.382578
$ws.Range("O10").Value = 0.441663298340486
$ws.Range("P10").Value = 0.441663298340486
$ws.Range("Q10").Value = 3.219127325120444
$ws.Range("R10").Value = 28.972145926084
$ws.Range("S10").Value = 0.1950664691009971
$ws.Range("T10").Value = 0.1950664691009971
